$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 25 and 26 (dimension shrinks to A1:F24)
$ws.Range("A25:A26").EntireRow.Delete()

# Row 2
$ws.Range("B2").Value = "T1"
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = "Mean (sd) : 0.3 (3.1)\`nmin < med < max:\`n0 < 0 < 99\`nIQR (CV) : 0 (10.5)"
$ws.Range("E2").Value = "0  : 1926 (91.1%)\`n1  :   53 ( 2.5%)\`n2  :   55 ( 2.6%)\`n3  :   48 ( 2.3%)\`n4  :   31 ( 1.5%)\`n99 :    2 ( 0.1%)"
$ws.Range("F2").Value = "0\`n(0%)"

# Row 3
$ws.Range("B3").Value = "T2"
$ws.Range("C3").ClearContents()
$ws.Range("D3").Value = "Mean (sd) : 1 (0.8)\`nmin < med < max:\`n0 < 1 < 2\`nIQR (CV) : 2 (0.8)"
$ws.Range("E3").Value = "0 : 669 (31.6%)\`n1 : 743 (35.1%)\`n2 : 703 (33.2%)"
$ws.Range("F3").Value = "0\`n(0%)"

# Row 4
$ws.Range("B4").Value = "T3"
$ws.Range("C4").ClearContents()
$ws.Range("D4").Value = "Min  : 0\`nMean : 0.1\`nMax  : 1"
$ws.Range("E4").Value = "0 : 1811 (85.6%)\`n1 :  304 (14.4%)"
$ws.Range("F4").Value = "0\`n(0%)"

# Row 5
$ws.Range("B5").Value = "T4"
$ws.Range("C5").ClearContents()
$ws.Range("D5").Value = "Mean (sd) : 1 (1.3)\`nmin < med < max:\`n0 < 1 < 4\`nIQR (CV) : 2 (1.3)"
$ws.Range("E5").Value = "0 : 1054 (49.8%)\`n1 :  491 (23.2%)\`n2 :  252 (11.9%)\`n3 :  132 ( 6.2%)\`n4 :  186 ( 8.8%)"
$ws.Range("F5").Value = "0\`n(0%)"

# Row 6
$ws.Range("B6").Value = "I1"
$ws.Range("C6").ClearContents()
$ws.Range("D6").Value = "Mean (sd) : 5.2 (16.5)\`nmin < med < max:\`n0 < 3 < 99\`nIQR (CV) : 1 (3.1)"
$ws.Range("E6").Value = "0  :   84 ( 4.0%)\`n1  :  231 (10.9%)\`n2  :  615 (29.1%)\`n3  : 1122 (53.0%)\`n99 :   63 ( 3.0%)"
$ws.Range("F6").Value = "0\`n(0%)"

# Row 7
$ws.Range("B7").Value = "I2"
$ws.Range("C7").ClearContents()
$ws.Range("D7").Value = "Mean (sd) : 1.5 (5.3)\`nmin < med < max:\`n0 < 1 < 99\`nIQR (CV) : 2 (3.5)"
$ws.Range("E7").Value = "0  : 580 (27.4%)\`n1  : 746 (35.3%)\`n2  : 489 (23.1%)\`n3  : 294 (13.9%)\`n99 :   6 ( 0.3%)"
$ws.Range("F7").Value = "0\`n(0%)"

# Row 8
$ws.Range("B8").Value = "I3"
$ws.Range("C8").ClearContents()
$ws.Range("D8").Value = "Mean (sd) : 2.5 (8.2)\`nmin < med < max:\`n0 < 2 < 99\`nIQR (CV) : 2 (3.2)"
$ws.Range("E8").Value = "0  : 295 (14.0%)\`n1  : 464 (21.9%)\`n2  : 592 (28.0%)\`n3  : 749 (35.4%)\`n99 :  15 ( 0.7%)"
$ws.Range("F8").Value = "0\`n(0%)"

# Row 9
$ws.Range("B9").Value = "D1"
$ws.Range("C9").ClearContents()
$ws.Range("D9").Value = "Mean (sd) : 0.9 (0.4)\`nmin < med < max:\`n0 < 1 < 2\`nIQR (CV) : 0 (0.5)"
$ws.Range("E9").Value = "0 :  351 (16.6%)\`n1 : 1715 (81.1%)\`n2 :   49 ( 2.3%)"
$ws.Range("F9").Value = "0\`n(0%)"

# Row 10
$ws.Range("B10").Value = "D2"
$ws.Range("C10").ClearContents()
$ws.Range("D10").Value = "Mean (sd) : 1 (0.3)\`nmin < med < max:\`n0 < 1 < 2\`nIQR (CV) : 0 (0.3)"
$ws.Range("E10").Value = "0 :   86 ( 4.1%)\`n1 : 1962 (92.8%)\`n2 :   67 ( 3.2%)"
$ws.Range("F10").Value = "0\`n(0%)"

# Row 11
$ws.Range("B11").Value = "D3"
$ws.Range("C11").ClearContents()
$ws.Range("D11").Value = "Mean (sd) : 0.2 (0.4)\`nmin < med < max:\`n0 < 0 < 2\`nIQR (CV) : 0 (2.4)"
$ws.Range("E11").Value = "0 : 1785 (84.4%)\`n1 :  295 (14.0%)\`n2 :   35 ( 1.6%)"
$ws.Range("F11").Value = "0\`n(0%)"

# Row 12
$ws.Range("B12").Value = "E1"
$ws.Range("C12").ClearContents()
$ws.Range("D12").Value = "Mean (sd) : 2.4 (12.4)\`nmin < med < max:\`n0 < 0 < 99\`nIQR (CV) : 2 (5.1)"
$ws.Range("E12").Value = "0  : 1222 (57.8%)\`n1  :  263 (12.4%)\`n2  :  252 (11.9%)\`n3  :  344 (16.3%)\`n99 :   34 ( 1.6%)"
$ws.Range("F12").Value = "0\`n(0%)"

# Row 13
$ws.Range("B13").Value = "E2"
$ws.Range("C13").ClearContents()
$ws.Range("D13").Value = "Mean (sd) : 4.9 (19.4)\`nmin < med < max:\`n0 < 0 < 99\`nIQR (CV) : 2 (3.9)"
$ws.Range("E13").Value = "0  : 1179 (55.7%)\`n1  :  206 ( 9.7%)\`n2  :  248 (11.7%)\`n3  :  396 (18.7%)\`n99 :   86 ( 4.1%)"
$ws.Range("F13").Value = "0\`n(0%)"

# Row 14
$ws.Range("B14").Value = "E3"
$ws.Range("C14").ClearContents()
$ws.Range("D14").Value = "Mean (sd) : 2.5 (13.6)\`nmin < med < max:\`n0 < 0 < 99\`nIQR (CV) : 1 (5.4)"
$ws.Range("E14").Value = "0  : 1417 (67.0%)\`n1  :  248 (11.7%)\`n2  :  217 (10.3%)\`n3  :  192 ( 9.1%)\`n99 :   41 ( 1.9%)"
$ws.Range("F14").Value = "0\`n(0%)"

# Row 15
$ws.Range("B15").Value = "E4"
$ws.Range("C15").ClearContents()
$ws.Range("D15").Value = "Mean (sd) : 3.7 (16.2)\`nmin < med < max:\`n0 < 0 < 99\`nIQR (CV) : 2 (4.4)"
$ws.Range("E15").Value = "0  : 1138 (53.8%)\`n1  :  235 (11.1%)\`n2  :  311 (14.7%)\`n3  :  372 (17.6%)\`n99 :   59 ( 2.8%)"
$ws.Range("F15").Value = "0\`n(0%)"

# Row 16
$ws.Range("B16").Value = "E5"
$ws.Range("C16").ClearContents()
$ws.Range("D16").Value = "Mean (sd) : 5.2 (20.6)\`nmin < med < max:\`n0 < 0 < 99\`nIQR (CV) : 2 (3.9)"
$ws.Range("E16").Value = "0  : 1308 (61.8%)\`n1  :  244 (11.5%)\`n2  :  212 (10.0%)\`n3  :  254 (12.0%)\`n99 :   97 ( 4.6%)"
$ws.Range("F16").Value = "0\`n(0%)"

# Row 17
$ws.Range("B17").Value = "CEX1"
$ws.Range("C17").ClearContents()
$ws.Range("D17").Value = "Min  : 0\`nMean : 0.1\`nMax  : 1"
$ws.Range("E17").Value = "0 : 1814 (85.8%)\`n1 :  301 (14.2%)"
$ws.Range("F17").Value = "0\`n(0%)"

# Row 18
$ws.Range("B18").Value = "CEX2"
$ws.Range("C18").ClearContents()
$ws.Range("D18").Value = "Mean (sd) : 0.5 (0.7)\`nmin < med < max:\`n0 < 0 < 2\`nIQR (CV) : 1 (1.3)"
$ws.Range("E18").Value = "0 : 1225 (57.9%)\`n1 :  633 (29.9%)\`n2 :  257 (12.2%)"
$ws.Range("F18").Value = "0\`n(0%)"

# Row 19
$ws.Range("B19").Value = "CEX3"
$ws.Range("C19").ClearContents()
$ws.Range("D19").Value = "Min  : 0\`nMean : 0.1\`nMax  : 1"
$ws.Range("E19").Value = "0 : 1842 (87.1%)\`n1 :  273 (12.9%)"
$ws.Range("F19").Value = "0\`n(0%)"

# Row 20
$ws.Range("B20").Value = "CEX4"
$ws.Range("C20").ClearContents()
$ws.Range("D20").Value = "Min  : 0\`nMean : 0.1\`nMax  : 1"
$ws.Range("E20").Value = "0 : 1968 (93.0%)\`n1 :  147 ( 7.0%)"
$ws.Range("F20").Value = "0\`n(0%)"

# Row 21
$ws.Range("B21").Value = "CEN1"
$ws.Range("C21").ClearContents()
$ws.Range("D21").Value = "Min  : 0\`nMean : 0.3\`nMax  : 1"
$ws.Range("E21").Value = "0 : 1379 (65.2%)\`n1 :  736 (34.8%)"
$ws.Range("F21").Value = "0\`n(0%)"

# Row 22
$ws.Range("B22").Value = "CEN2"
$ws.Range("C22").ClearContents()
$ws.Range("D22").Value = "Min  : 0\`nMean : 0.8\`nMax  : 1"
$ws.Range("E22").Value = "0 :  376 (17.8%)\`n1 : 1739 (82.2%)"
$ws.Range("F22").Value = "0\`n(0%)"

# Row 23
$ws.Range("B23").Value = "CEN3"
$ws.Range("C23").ClearContents()
$ws.Range("D23").Value = "Min  : 0\`nMean : 0.5\`nMax  : 1"
$ws.Range("E23").Value = "0 : 1140 (53.9%)\`n1 :  975 (46.1%)"
$ws.Range("F23").Value = "0\`n(0%)"

# Row 24
$ws.Range("B24").Value = "FEXP_MAY15"
$ws.Range("C24").Value = "FACTOR EXPANSIÓN"
$ws.Range("D24").Value = "Mean (sd) : 726.6 (1455.6)\`nmin < med < max:\`n4.8 < 307.1 < 20920.1\`nIQR (CV) : 628.1 (2)"
$ws.Range("E24").Value = "1747 distinct values"
$ws.Range("F24").Value = "0\`n(0%)"

